$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# "Description"/"Paths" columns from A/B to B/C.
$ws.Columns("A:A").Insert()

# Fill in the new "Title" column.
$ws.Range("A1").Value = "Title"
$ws.Range("A2").Value = "Film Suggestion App"

# Update the Description text (now in column B) to mention "streamlit".
$ws.Range("B2").Value = "It's a streamlit web app that suggests films from a category that user's selected"

# Restore / set the column widths to match the final layout.
$ws.Columns("A").ColumnWidth = 57
$ws.Columns("B").ColumnWidth = 66
$ws.Columns("C").ColumnWidth = 44.666666666666664

# Match the final selection state.
$ws.Range("A3").Select() | Out-Null
